$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay stored as text
# (matching the source data which is all inline/shared strings),
# so force a Text number format before writing the value.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply the updated coin data
$ws.Range('D2').Value = '51.206.03'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '2.950.96'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '380.54'
$ws.Range('E5').Value = '  +2.65%  '
$ws.Range('D6').Value = '104.61'
$ws.Range('E6').Value = '  +1.36%  '
$ws.Range('D7').Value = '0.539'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('D10').Value = '36.93'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  +0.67%  '
$ws.Range('D12').Value = '0.0839'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '3.420.46'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('D14').Value = '18.37'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '7.45'
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '2.962.93'
$ws.Range('E16').Value = '  +1.37%  '
$ws.Range('D17').Value = '0.957'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').Value = '51.161.29'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '3.33'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').Value = '7.33'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '12.85'
$ws.Range('E21').Value = '  -1.25%  '
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('D23').Value = '68.84'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('D24').Value = '260.38'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '2.81'
$ws.Range('E25').Value = '  +3.62%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.168'
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('B27').Value = 'Filecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D27').Value = '7.45'
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('D28').Value = '7.12'
$ws.Range('E28').Value = '  +16.47%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '25.82'
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  +8.16%  '
$ws.Range('D32').Value = '9.80'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').Value = '34.53'
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('D35').Value = '51.06'
$ws.Range('E35').Value = '  +1.00%  '
$ws.Range('D36').Value = '0.0446'
$ws.Range('E36').Value = '  +5.59%  '
$ws.Range('D37').Value = '1.01'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '3.05'
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('D39').Value = '17.23'
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('D40').Value = '2.57'
$ws.Range('E40').Value = '  -4.23%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('E42').Value = '  +1.82%  '
$ws.Range('D43').Value = '122.75'
$ws.Range('E43').Value = '  +4.12%  '
$ws.Range('D44').Value = '21.83'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').Value = '0.279'
$ws.Range('E45').Value = '  +16.23%  '
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').Value = '2.033.20'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('D49').Value = '3.20'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('D50').Value = '0.0348'
$ws.Range('E50').Value = '  +10.08%  '
$ws.Range('D51').Value = '5.10'
$ws.Range('E51').Value = '  +1.35%  '
